$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

$cellTypes = $sheets.Item("cell_types")
$immuneCells = $sheets.Item("immune_cells")

# Insert the new sheet right before "immune_cells" (i.e. right after "cell_types")
$new = $sheets.Add($null, $cellTypes)
$new.Name = "cell_types_for_tcga"

$data = @(
    @("value", "color"),
    @("T cell regulatory (Tregs)", "#8dd3c7"),
    @("T cell CD4+", "#8dd3c7"),
    @("T cell CD8+", "#ffffb3"),
    @("NK cell", "#ffed6f"),
    @("Dendritic cell", "#fb8072"),
    @("Neutrophil", "#80b1d3"),
    @("Monocyte", "#fdb462"),
    @("Macrophage", "#fdb462"),
    @("Macrophage/Monocyte", "#fdb462"),
    @("B cell", "#b3de69"),
    @("Endothelial cell", "#fccde5"),
    @("Cancer associated fibroblast", "#bc80bd"),
    @("Eosinophil", "#999999")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $new.Cells.Item($row, 1).Value = $data[$i][0]
    $new.Cells.Item($row, 2).Value = $data[$i][1]
}

$new.Columns.Item(1).ColumnWidth = 27.140625

# The previously-active "immune_cells" sheet loses its old selection/activation.
# (Re-fetch by name: the old $immuneCells handle now resolves to the newly
# inserted sheet since worksheet handles track by position, not identity.)
$immuneCellsAfter = $sheets.Item("immune_cells")
$immuneCellsAfter.Range("A1:B10").Select() | Out-Null

# Select the new sheet last so it ends up as the active tab, matching the
# saved selection state for "cell_types_for_tcga".
$new.Range("H17").Select() | Out-Null
